# Rename all worksheets in the workbook according to the new naming scheme.
# The sheets keep their relative order (sheetId / r:id / position unchanged),
# only the visible sheet name changes.

$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ16786044",
    "summ16951479",
    "summ17252913",
    "summ17548580",
    "summ17816064",
    "summ18047654",
    "summ18268288",
    "summ18481092",
    "summ18677793",
    "summ18881440",
    "summ19103381",
    "summ19323430",
    "summ19530227",
    "summ19739905",
    "summ19945285",
    "summ20191814",
    "summ20442261",
    "summ20672624",
    "summ20897768",
    "summ21190819",
    "summ21399143",
    "summ21590071",
    "summ21878270",
    "summ22104355",
    "summ22350477",
    "summ22563761",
    "summ22771678",
    "summ23002505",
    "summ23232331",
    "summ23484876",
    "summ23705524",
    "summ23916865",
    "summ24116888",
    "summ24312977",
    "summ24540603",
    "summ24751844",
    "summ24991060",
    "summ25242123",
    "summ25512795",
    "summ25793450",
    "summ26026860",
    "summ26243234",
    "summ26467892",
    "summ26677800",
    "summ26925437",
    "summ27173691",
    "summ27490259",
    "summ27732421",
    "summ27973326",
    "summ28206256"
)

for ($i = 1; $i -le $newNames.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]
}
